$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row gets a third column: "Role"
$ws.Range("C1").Value = "Role"

# Existing login row (row 2) now also carries a Role of "Admin"
$ws.Range("C2").Value = "Admin"

# New rows 3-10: same Username ("ATPBoss"), a sequence of new failed
# passwords, and the same Role ("Admin") - additional negative test data
# for the separated "Excel Reader" test case.
$ws.Range("A3").Value = "ATPBoss"
$ws.Range("B3").Value = "Fail1235"
$ws.Range("C3").Value = "Admin"

$ws.Range("A4").Value = "ATPBoss"
$ws.Range("B4").Value = "Fail1236"
$ws.Range("C4").Value = "Admin"

$ws.Range("A5").Value = "ATPBoss"
$ws.Range("B5").Value = "Fail1237"
$ws.Range("C5").Value = "Admin"

$ws.Range("A6").Value = "ATPBoss"
$ws.Range("B6").Value = "Fail1238"
$ws.Range("C6").Value = "Admin"

$ws.Range("A7").Value = "ATPBoss"
$ws.Range("B7").Value = "Fail1239"
$ws.Range("C7").Value = "Admin"

$ws.Range("A8").Value = "ATPBoss"
$ws.Range("B8").Value = "Fail1240"
$ws.Range("C8").Value = "Admin"

$ws.Range("A9").Value = "ATPBoss"
$ws.Range("B9").Value = "Fail1241"
$ws.Range("C9").Value = "Admin"

$ws.Range("A10").Value = "ATPBoss"
$ws.Range("B10").Value = "Fail1242"
$ws.Range("C10").Value = "Admin"

# Match the authored file's recorded active-cell selection
$ws.Range("H17").Select()
